$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'criterion'
$ws.Range("B1").Value = 'mercado_libre'
$ws.Range("C1").Value = 'amazon'
$ws.Range("D1").Value = 'pedidoya'
$ws.Range("E1").Value = 'website_from_image'
$ws.Range("F1").Value = 'conclusion'
$ws.Range("A2").Value = 'Typography'
$ws.Range("B2").Value = 'Utiliza una tipografía Sans-serif clara y legible, como Arial o Open Sans, que asegura una buena experiencia de lectura en todos sus dispositivos.'
$ws.Range("C2").Value = 'Emplea una tipografía distintiva (Amazon Ember) para mantener consistencia visual, facilitando la identificación de la marca y la lectura de descripciones de productos.'
$ws.Range("D2").Value = 'Su tipografía es moderna y amigable, a menudo usando una fuente redondeada y sencilla para transmitir accesibilidad y rapidez en el servicio de entrega.'
$ws.Range("E2").Value = 'La tipografía es limpia y moderna, principalmente sans-serif, ofreciendo buena legibilidad para nombres de productos y encabezados importantes en la página.'
$ws.Range("F2").Value = 'Considerar incorporar variaciones de peso o tamaño para destacar mejor elementos clave sin comprometer la consistencia global del sitio.'
$ws.Range("A3").Value = 'Colors'
$ws.Range("B3").Value = 'Domina el amarillo vibrante junto con el azul oscuro, creando un contraste energético que llama la atención y resalta ofertas especiales.'
$ws.Range("C3").Value = 'Predominan el blanco, negro y un toque de naranja/azul, generando un esquema de colores profesional que facilita la concentración en los productos mostrados.'
$ws.Range("D3").Value = 'Se caracteriza por el rojo brillante y el blanco, transmitiendo energía y urgencia, lo cual es muy apropiado para su modelo de negocio de entrega rápida.'
$ws.Range("E3").Value = 'La paleta se centra en un azul brillante para la cabecera y el blanco para el contenido, generando una apariencia limpia y enfocada en los productos.'
$ws.Range("F3").Value = 'Explorar el uso de tonos complementarios sutiles para áreas interactivas, mejorando la jerarquía visual sin saturar el diseño existente.'
$ws.Range("B4").Value = 'Adopta un tono semi-formal, equilibrando la profesionalidad de las transacciones con un lenguaje cercano para construir confianza con sus usuarios.'
$ws.Range("C4").Value = 'Mantiene un tono predominantemente formal y transaccional, priorizando la claridad y eficiencia en la comunicación de detalles de productos y servicios.'
$ws.Range("D4").Value = 'Utiliza un lenguaje muy informal y amigable, con emojis y expresiones coloquiales que reflejan la juventud y dinamismo de su público objetivo.'
$ws.Range("E4").Value = 'Muestra un tono más formal y directo, enfocado en la funcionalidad y la información técnica, adecuado para un nicho de componentes tecnológicos.'
$ws.Range("F4").Value = 'Podría beneficiarse de introducir elementos de lenguaje ligeramente más conversacionales para fomentar mayor engagement con la comunidad de ensambladores.'
$ws.Range("A5").Value = 'Characters / Icons / Emblems'
$ws.Range("B5").Value = 'Su logo con la balanza y flecha simboliza comercio y eficiencia; los iconos son minimalistas y funcionales, facilitando la navegación sin distracciones.'
$ws.Range("C5").Value = 'El logo con la flecha de la ''a'' a la ''z'' subraya su vasta oferta; los iconos son universales y muy reconocibles, promoviendo una interacción intuitiva.'
$ws.Range("D5").Value = 'Usa un logo sencillo con un globo de diálogo y un tenedor, junto a iconos claros que representan categorías de comida y estados de pedidos.'
$ws.Range("E5").Value = 'El logo es una combinación de un ratón de PC con un embudo/circuito, y los iconos (lupa, persona) son estándar y claros, mejorando la usabilidad.'
$ws.Range("F5").Value = 'Desarrollar un conjunto de iconos personalizados y únicos que reflejen aún más la temática de construcción de PCs, añadiendo personalidad de marca.'
$ws.Range("A6").Value = 'Accessibility'
$ws.Range("B6").Value = 'Ofrece opciones como descripciones de imágenes y compatibilidad con lectores de pantalla, buscando ser inclusivo para usuarios con diversas capacidades.'
$ws.Range("C6").Value = 'Invierte en características de accesibilidad robustas, incluyendo navegación por teclado y soporte para texto de alto contraste, garantizando una experiencia para todos.'
$ws.Range("D6").Value = 'Se enfoca en una interfaz sencilla con grandes botones, facilitando el uso para usuarios con posibles dificultades motoras o visuales leves, priorizando la rapidez.'
$ws.Range("E6").Value = 'La legibilidad del texto y el contraste de colores son adecuados, lo que sugiere una buena base para la accesibilidad general de la plataforma.'
$ws.Range("F6").Value = 'Implementar alternativas de texto para todas las imágenes y asegurar una navegación totalmente controlable por teclado para usuarios con discapacidades.'
$ws.Range("A7").Value = 'Navigation (important buttons)'
$ws.Range("B7").Value = 'Botones de "Comprar" y "Agregar al carrito" son prominentes, junto a una barra de búsqueda eficiente, guiando al usuario rápidamente hacia la compra.'
$ws.Range("C7").Value = 'Destaca por su navegación intuitiva con un mega-menú y botones de acción claros, permitiendo a los usuarios encontrar productos con facilidad.'
$ws.Range("D7").Value = 'Sus botones de "Pedir" y "Reordenar" son muy visibles, con categorías de restaurantes claras, facilitando la elección y confirmación rápida de pedidos.'
$ws.Range("E7").Value = 'Los botones "Información", "Arma tu PC", "Comparar" son claros en la cabecera, junto a un campo de búsqueda prominente para encontrar productos.'
$ws.Range("F7").Value = 'Mejorar la visibilidad de los botones de filtros y opciones de ordenamiento dentro de la sección de productos para una mejor exploración.'
$ws.Range("A8").Value = 'Organization'
$ws.Range("B8").Value = 'La organización de productos por categorías, tiendas oficiales y secciones de ofertas es clara, permitiendo una búsqueda eficiente para los usuarios.'
$ws.Range("C8").Value = 'Utiliza una categorización detallada y filtros avanzados, organizando eficazmente millones de productos para una experiencia de compra personalizada y sencilla.'
$ws.Range("D8").Value = 'La organización por tipo de cocina, promociones y distancia es eficiente, permitiendo a los usuarios descubrir y seleccionar restaurantes rápidamente.'
$ws.Range("E8").Value = 'Presenta los "Componentes populares" de manera atractiva en cuadrícula, con títulos de productos claros, mostrando una organización lógica inicial.'
$ws.Range("F8").Value = 'Incluir más opciones de filtrado y ordenamiento detallado (por socket, generación, marca) para los componentes populares mejoraría la experiencia.'
$ws.Range("A9").Value = 'Extra features'
$ws.Range("B9").Value = 'Ofrece Mercado Puntos, envíos Flex, y un sistema de preguntas y respuestas, añadiendo valor a la experiencia de compra y fidelización del cliente.'
$ws.Range("C9").Value = 'Presenta Amazon Prime, reseñas detalladas, listas de deseos y recomendaciones personalizadas, enriqueciendo significativamente la interacción del usuario.'
$ws.Range("D9").Value = 'Dispone de seguimiento de pedidos en tiempo real, promociones exclusivas y la opción de calificar restaurantes, mejorando la conveniencia del servicio.'
$ws.Range("E9").Value = 'Ofrece "Arma tu PC" y "Comparar" como características distintivas, lo cual es muy útil para su audiencia especializada en hardware de computadoras.'
$ws.Range("F9").Value = 'Considerar la implementación de un foro comunitario o un blog con guías de ensamblaje para fomentar la interacción y educar a los usuarios.'
$ws.Range("A10").Value = 'Tutorials or Instructions'
$ws.Range("B10").Value = 'Proporciona secciones de ayuda y guías para vendedores y compradores, facilitando el uso de la plataforma y resolviendo dudas comunes.'
$ws.Range("C10").Value = 'Cuenta con amplias páginas de ayuda, tutoriales en video y descripciones detalladas de productos, asistiendo a los usuarios en cada etapa de su compra.'
$ws.Range("D10").Value = 'Sus instrucciones son simples y directas para realizar pedidos, con una sección de preguntas frecuentes que resuelve dudas básicas de manera eficiente.'
$ws.Range("E10").Value = 'La sección "Arma tu PC" sugiere que hay guías implícitas, pero una sección explícita de "Información" o "Ayuda" sería beneficiosa y más visible.'
$ws.Range("F10").Value = 'Desarrollar una base de conocimientos completa con tutoriales detallados sobre la compatibilidad de componentes y guías de ensamblaje paso a paso.'
$ws.Range("A11").Value = 'Overall User Experience'
$ws.Range("B11").Value = 'Ofrece una experiencia completa y funcional, con facilidad de compra y venta, aunque a veces la saturación de ofertas puede ser abrumadora.'
$ws.Range("C11").Value = 'Proporciona una UX pulida y eficiente, priorizando la comodidad y personalización, lo que permite a los usuarios encontrar y comprar con facilidad.'
$ws.Range("D11").Value = 'Se enfoca en una experiencia rápida y sin fricciones para pedir comida, con una interfaz intuitiva que minimiza los pasos necesarios para completar la orden.'
$ws.Range("E11").Value = 'La interfaz es limpia y enfocada, facilitando la visualización de componentes y la navegación básica para usuarios interesados en ensamblar PCs.'
$ws.Range("F11").Value = 'Optimizar los tiempos de carga de las páginas de productos y mejorar la retroalimentación visual en las interacciones para una experiencia más fluida.'
